# Rotate the species-observation records currently on rows 3, 4 and 5:
#   new row 3 <- old row 4
#   new row 4 <- old row 5
#   new row 5 <- old row 3
# (all other rows / columns are left untouched)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstCol = "A"
$lastCol  = "AY"
$tempRow  = 100   # scratch row far away from the real data, cleared at the end

$srcRow3 = "$firstCol" + "3:" + "$lastCol" + "3"
$srcRow4 = "$firstCol" + "4:" + "$lastCol" + "4"
$srcRow5 = "$firstCol" + "5:" + "$lastCol" + "5"
$tempRng = "$firstCol" + "$tempRow" + ":" + "$lastCol" + "$tempRow"

# Stash the current row 3 in a scratch row (Copy preserves value types,
# e.g. text-that-looks-numeric stays text instead of being re-parsed).
$ws.Range($tempRng).ClearContents()
$ws.Range($srcRow3).Copy($ws.Range($tempRng))

# row3 = old row4  (clear first: Copy() only overwrites cells that are
# non-empty in the source, so stale values must be wiped first)
$ws.Range($srcRow3).ClearContents()
$ws.Range($srcRow4).Copy($ws.Range($srcRow3))

# row4 = old row5
$ws.Range($srcRow4).ClearContents()
$ws.Range($srcRow5).Copy($ws.Range($srcRow4))

# row5 = old row3 (retrieved from the scratch row)
$ws.Range($srcRow5).ClearContents()
$ws.Range($tempRng).Copy($ws.Range($srcRow5))

# tidy up the scratch row
$ws.Range($tempRng).ClearContents()
